# Update "想去人数" (number of people interested) figures on the
# 展览 (Exhibitions) and 全部类型 (All types) sheets to match the
# freshly generated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 113
$wsExpo.Range("F5").Value = 17
$wsExpo.Range("F6").Value = 0
$wsExpo.Range("F7").Value = 0
$wsExpo.Range("F8").Value = 143
$wsExpo.Range("F9").Value = 63
$wsExpo.Range("F10").Value = 0

# --- Sheet: 全部类型 -------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 113
$wsAll.Range("F4").Value = 1614
$wsAll.Range("F5").Value = 17
$wsAll.Range("F8").Value = 143
$wsAll.Range("F10").Value = 478
